$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial date number): 45406 -> 45436
$ws.Range("A1").Value = 45436

# Update the price values in column D for rows 31-33
$ws.Range("D31").Value = 3985.166
$ws.Range("D32").Value = 1646.582
$ws.Range("D33").Value = 3379.108
